$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 2).Value = "ce8ab73d0b5eea84312e9ae5aaec7d18"  # 05-050305TC (was 6c1a51b5e85289c72b553ad899db1a51)
$ws.Cells.Item(11, 2).Value = "1352d9b99bf06626ff80952eda02d7d2"  # 05-050301A (was 1f682c4baf00039722b9d3b2a8f6431f)
$ws.Cells.Item(15, 2).Value = "454bd2ded7f4437defc0cd535776b96c"  # 05-050207TP (was 748fdfa86f292b105e0f90f30045b1c7)
$ws.Cells.Item(17, 2).Value = "c41c5b03bfac8c0e8061ad380b8aee67"  # 05-050305TP (was 9f0009d40270a79205ceb4b19fbe61e3)
$ws.Cells.Item(24, 2).Value = "5cfb543c596d343f616d0f935a2d1c0c"  # 05-050316TC (was d6e0c50a94994e93363908ad3893b5fb)
$ws.Cells.Item(29, 2).Value = "efede8d45f30e4682811f2ed6b25fea8"  # 05-050302A (was 14bcaddadb80968b9f4d54b9a0bf5b4d)
$ws.Cells.Item(34, 2).Value = "c61e0c5fa0c3d3aeb7f195c62229f494"  # 05-050316TP (was 9b5fa738b68a8c46f512c3e8ae609d3b)
$ws.Cells.Item(100, 2).Value = "aed0b10bec2e9f418199ec1ba1e7d195"  # 04-040021TM (was 85819c9b0ee864700a6fb3abf7b62758)
$ws.Cells.Item(104, 2).Value = "d3250a5c1d6eed791df60eeb4e0dcd1e"  # 04-040021TP (was afc45b0ea45fcd2114d8102997488408)
$ws.Cells.Item(113, 2).Value = "d2c1e61c41dbe7d42161673c62f26d94"  # 04-040021TC (was 956b266fd844e9f3fca2194ee278fadb)
$ws.Cells.Item(121, 2).Value = "27c1bb70cb640d5ca20a759347c927c8"  # 05-050301TP (was 81667d4f5140992663fc6287a415e11f)
$ws.Cells.Item(122, 2).Value = "6403db4eaca423e88668dcf035f15b05"  # 04-040014TC (was d15ca3c8fb72fbbd22db7c2394f28a69)
$ws.Cells.Item(133, 2).Value = "c93f48efa2270bb475f8fe1ae270e4eb"  # 05-050312TP (was e67cb7acf6fa4ff9ebae00920bc5988a)
$ws.Cells.Item(136, 2).Value = "9d3d20be9dae4a3465495cdc14ee1b80"  # 05-050312TC (was 4f5900008902af644f9584451a1c3901)
$ws.Cells.Item(159, 2).Value = "ab8b7ee5396a905857d3beda3e826fb5"  # 05-050203TP (was 2300b18945809dfc46c117b49b348318)
$ws.Cells.Item(162, 2).Value = "28b7081ddd8b2bf574091a34d8703cef"  # 05-050308A (was 537a5222143850acb0b8e7c2a56d1a6f)
$ws.Cells.Item(164, 2).Value = "2c5adae7a14dbf122a3e7e333946951c"  # 04-040021A (was 0a80cf60deec27272e68c8141fbee685)
$ws.Cells.Item(169, 2).Value = "7196fce0dc6e111edfa1e0070365cd31"  # 05-050203TC (was 3f3b4e700fa10a31b1e86ddb99367c7c)
$ws.Cells.Item(180, 2).Value = "4452182d4a3e39871668d09fdb6c1e5b"  # 05-050303TC (was 8e3e66726412138b9c21d57bc4009d98)
$ws.Cells.Item(183, 2).Value = "a52b70165fb48df54710bb99294990ad"  # 05-050305A (was b2ea477540860dd093eec216119c4402)
$ws.Cells.Item(191, 2).Value = "c73e5ad0a567948972aa3db3a087d497"  # 05-050314TP (was 6aeb8c7ff9726e431785864e912f5be0)
$ws.Cells.Item(198, 2).Value = "d4be0e2477c0204b63749618f25577c3"  # 05-050314TC (was d3bf3c1c93e8e11b73485fcb6846fff5)
$ws.Cells.Item(200, 2).Value = "a40b5de7c55ab8e14ffb3a6cb9a21819"  # 05-050306A (was b123bd7cd912a41db92fb5ee74f564a3)
$ws.Cells.Item(213, 2).Value = "e11742ebab986b101aaf472dd8371e81"  # 05-050303A (was f1a3da6a4991d211f4d0e18b9486ed7a)
$ws.Cells.Item(228, 2).Value = "999429de30030d0246f871e7b12f1894"  # 05-050304A (was da137e8bd5d8f137f86514581a664b40)
$ws.Cells.Item(230, 2).Value = "6cc38a03a89a547d65027e64ed10d11e"  # 04-040014TP (was a7ccd9496d18261177551264266f67e7)
$ws.Cells.Item(233, 2).Value = "bf8a8bb894e8abfbc617dca6f524925d"  # 04-040014TM (was 380c5e4c6ed05e85df43317f9a0cfa66)
$ws.Cells.Item(281, 2).Value = "128753153d13aec6cc2385e89c1d0dfa"  # 05-050201TC (was 0d957c92e3ee8329abd131457daee2a1)
$ws.Cells.Item(331, 2).Value = "e40c86b9f34ec9b35c007636ab13d92c"  # 04-040018TP (was d9986ed4380897b50d61c0803314de7c)
$ws.Cells.Item(339, 2).Value = "ee3d156fa6c8abe62caff64a715ba9ef"  # 05-050201TP (was 391520f32aea89a505c0189d4b31d852)
$ws.Cells.Item(342, 2).Value = "091253ce18fb2eaae2c3e52d7191a868"  # 04-040018TC (was 052d5b4453144717d9154004c40aed09)
$ws.Cells.Item(343, 2).Value = "9683867abb9ed08c897898b9ce16f688"  # 04-040018TM (was 9c8e173b79f48d63f00af95644862e76)
$ws.Cells.Item(381, 2).Value = "ccb51bd55ef71d785c4cbe725d27c184"  # 01-010073A (was 426758b07b194188b97fe09b886f440d)
$ws.Cells.Item(419, 2).Value = "bf3569543f5afe0bd329968445d710df"  # 05-0709-070905BTC (was 930e9bd628ccd09c643cd2b4a4b8cfad)
$ws.Cells.Item(458, 2).Value = "752988414c894035dd2770010236af04"  # 01-010073TP (was 62f05aaa5756711c583f9c74bdffd409)
$ws.Cells.Item(461, 2).Value = "b11b80ec3b93464d6b97a5f9c1948435"  # 05-050313A (was 060072cb4a449d58d07838c00b609f70)
$ws.Cells.Item(477, 2).Value = "d42521fa4802f5f3088dfd72d207e8c7"  # 01-010073TC (was e1b8840a7130774ea1c4a2335241f85b)
$ws.Cells.Item(480, 2).Value = "f23b3dca7b162c63f81a3379142179f4"  # 05-050314A (was c2cefcf8311326ea2d84c3e9ddd5d4ad)
$ws.Cells.Item(502, 2).Value = "7a7522b01202c942ed0d664c6fa3c80c"  # 05-050208TP (was df800795b07776270ef538ef1fc07b06)
$ws.Cells.Item(506, 2).Value = "51d94fbb108c060af0774f3dfc25fd2e"  # 05-050306TP (was aa1791820592e49d2dde3aff5748084a)
$ws.Cells.Item(514, 2).Value = "1ff4dd27e25e4cecffa8c888a063c5c2"  # 05-050317TC (was 0163ad4ebad868ebcb1fb1d515410e6b)
$ws.Cells.Item(524, 2).Value = "586802b4d9ba45de50d961c63708f3c0"  # 05-050317TP (was b8463e643f40c14c051b7aa3e19cc647)
$ws.Cells.Item(563, 2).Value = "df2450d26af44f1ec23f8f2aca1c0b8f"  # 05-050308TC (was 58aeeda8ebd6873d630280821cb636b9)
$ws.Cells.Item(572, 2).Value = "f1eff8d1240251c266d684e4cbc1fca7"  # 05-050308TP (was 5ed55f8b2ae0bd9cea467720286f267b)
$ws.Cells.Item(616, 2).Value = "858655ce5cd775efbf95ef913749a8e1"  # 05-050204TP (was 205570524adee08761635f152af1eadb)
$ws.Cells.Item(619, 2).Value = "5f08ecffecd63e81c0870b802b54b76b"  # 04-040015TC (was bd09cfb4e9f5a5a1edc58ee2f6cbef23)
$ws.Cells.Item(623, 2).Value = "ff9f888e91bca8d85efafc7661513a32"  # 04-040015TP (was 5df9e1ffb7ca51b90d6720532ccfee6f)
$ws.Cells.Item(627, 2).Value = "366b49e650fa84c6e5ce2262d5375666"  # 05-050204TC (was b75d6c03eda2947c2f583a157129c161)
$ws.Cells.Item(628, 2).Value = "846c9647ded4ae397a5a92e7ec1d0301"  # 04-040015TM (was ae8a27b09551a4de674da30e82a0e23c)
$ws.Cells.Item(629, 2).Value = "326ca7636a2692fab909e88061b27250"  # 05-050302TP (was b4bf40be839e72ff90e5a588136c4567)
$ws.Cells.Item(666, 2).Value = "6a504f8d367e29df8fe91b6e061f2350"  # 05-050317A (was d0198b482e7ad0701fea272aba6657a8)
$ws.Cells.Item(680, 2).Value = "284db2d061666b0408f17f55094b0c4b"  # 05-050206TP (was 902b8c6b60528c7b830052fb906e50a8)
$ws.Cells.Item(685, 2).Value = "422a1417ff76cacbd629ea326069a450"  # 05-050206TC (was 225498260d678337a4782766a1156652)
$ws.Cells.Item(700, 2).Value = "c53e1c9f156064338c17d271449db459"  # 05-050304TC (was 54f8f0d13d2be919db718fbd6002f7de)
$ws.Cells.Item(703, 2).Value = "d216ae7caa824e3fcb0fa8c834c19559"  # 05-050206A (was 09d87cbc478370a8a2f110e3b1786283)
$ws.Cells.Item(704, 2).Value = "dfafc1925e9040bbb89dabf3f5bcc796"  # 05-050315A (was 0c15d29fc30a8a3b76d70a057ca66b27)
$ws.Cells.Item(715, 2).Value = "882e675b19e36e77fa97af68b7bce65f"  # 05-050304TP (was 6fec891a7daf86028b2467a7fac67a3f)
$ws.Cells.Item(729, 2).Value = "27ed38bf1fbffac7273df8279ccba7ca"  # 05-050316A (was b4db0bd5cfe9f51ea71702c7935a8b82)
$ws.Cells.Item(742, 2).Value = "a635002fc3375c746f355b5ba251ce4c"  # 05-050315TP (was 3945cc1ced32bc3ccd9b183feb1b5bcd)
$ws.Cells.Item(779, 2).Value = "d7f4356c35eb2b8b0deaac7d4e0be00c"  # 04-040018A (was babf3fd530aff2ea45435a4292853ff1)
$ws.Cells.Item(818, 2).Value = "1dcbd17e31672161575ab11d3dad0626"  # 04-040015A (was 4c2ed9e49577e877cba8646fab52dc00)
$ws.Cells.Item(819, 2).Value = "a729ea1a0ad69ec9de8500cc84a9dde0"  # 05-050202TP (was 34118f945a133ca7e014c23ed75edc68)
$ws.Cells.Item(830, 2).Value = "47aa2897ed7808fd9551156aea068daf"  # 05-050311TC (was f1a61ae09a06993f94701cb2daa2fa3d)
$ws.Cells.Item(831, 2).Value = "572bd04638f6b5cead7fe5e2de230d72"  # 04-040014A (was 3ebef27ff7385eb5bb0c6c1d9dc07834)
$ws.Cells.Item(835, 2).Value = "485754eddf6db83a63940e6505b91915"  # 05-050311TP (was 493485141f8ff34952434469c68d6932)
$ws.Cells.Item(874, 2).Value = "d878f735a89572d2273c1e98708e28dd"  # 03-030032A (was c9c849f03081bb7a17b5eba5feebb7ea)
